# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") is re-derived from the refreshed source data; write
# the newly computed K values for each data row (rows 2-46) in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 3
    8  = 3
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 3
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 0
    21 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 2
    31 = 3
    32 = 1
    33 = 2
    34 = 0
    35 = 1
    36 = 2
    37 = 3
    38 = 2
    39 = 1
    40 = 3
    41 = 2
    42 = 4
    43 = 0
    44 = 1
    45 = 1
    46 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
